$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fill in the "Java Default Sort" timing column (I) with the measured
# trial results / "Unmeasureable" markers that mirror column H's pattern
# for each input-size block.

$ws.Range("I6").Value = "Unmeasureable"
$ws.Range("I7").Value = "Unmeasureable"
$ws.Range("I8").Value = "Unmeasureable"
$ws.Range("I9").Value = "Unmeasureable"
$ws.Range("I10").Value = "Unmeasureable"

$ws.Range("I11").Value = 3521
$ws.Range("I12").Value = 3527
$ws.Range("I13").Value = 3546
$ws.Range("I14").Value = 3585
$ws.Range("I15").Value = 3555

$ws.Range("I16").Value = "Unmeasureable"
$ws.Range("I17").Value = "Unmeasureable"
$ws.Range("I18").Value = "Unmeasureable"
$ws.Range("I19").Value = "Unmeasureable"
$ws.Range("I20").Value = "Unmeasureable"

$ws.Range("I21").Value = 2597
$ws.Range("I22").Value = 2607
$ws.Range("I23").Value = 1666
$ws.Range("I24").Value = 1695
$ws.Range("I25").Value = 1676

$ws.Range("I26").Value = "Unmeasureable"
$ws.Range("I27").Value = "Unmeasureable"
$ws.Range("I28").Value = "Unmeasureable"
$ws.Range("I29").Value = "Unmeasureable"
$ws.Range("I30").Value = "Unmeasureable"

$ws.Range("I31").Value = 3373
$ws.Range("I32").Value = 3296
$ws.Range("I33").Value = 3280
$ws.Range("I34").Value = 3287
$ws.Range("I35").Value = 3402

$ws.Range("I36").Value = 1441
$ws.Range("I37").Value = 1332
$ws.Range("I38").Value = 1322
$ws.Range("I39").Value = 1350
$ws.Range("I40").Value = 1319

$ws.Range("I41").Value = 3533
$ws.Range("I42").Value = 3242
$ws.Range("I43").Value = 3245
$ws.Range("I44").Value = 3263
$ws.Range("I45").Value = 3259

$ws.Range("I46").Value = 1232
$ws.Range("I47").Value = 1165
$ws.Range("I48").Value = 1178
$ws.Range("I49").Value = 1176
$ws.Range("I50").Value = 1185

$ws.Range("I51").Value = 65
$ws.Range("I52").Value = 23
$ws.Range("I53").Value = 24
$ws.Range("I54").Value = 18
$ws.Range("I55").Value = 17

$ws.Range("I56").Value = 26
$ws.Range("I57").Value = 16
$ws.Range("I58").Value = 14
$ws.Range("I59").Value = 12
$ws.Range("I60").Value = 13

$ws.Range("I61").Value = 16
$ws.Range("I62").Value = 11
$ws.Range("I63").Value = 8
$ws.Range("I64").Value = 8
$ws.Range("I65").Value = 9

$ws.Range("I66").Value = 1387
$ws.Range("I67").Value = 1183
$ws.Range("I68").Value = 1344
$ws.Range("I69").Value = 1384
$ws.Range("I70").Value = 1373

$ws.Range("I71").Value = 319
$ws.Range("I72").Value = 310
$ws.Range("I73").Value = 350
$ws.Range("I74").Value = 306
$ws.Range("I75").Value = 303

$ws.Range("I76").Value = 1438
$ws.Range("I77").Value = 1419
$ws.Range("I78").Value = 1445
$ws.Range("I79").Value = 1416
$ws.Range("I80").Value = 1415

$ws.Range("I81").Value = 322
$ws.Range("I82").Value = 359
$ws.Range("I83").Value = 300
$ws.Range("I84").Value = 301
$ws.Range("I85").Value = 307

# Match the saved view state: scrolled down a bit further, with the new
# last cell of data (I86) now the active selection.
$ws.Range("I86").Select()
